$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 74 (shifts existing rows 74-83 down to 75-84)
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(74, 1).Value = 5
$ws.Cells.Item(74, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value = "Maule"
$ws.Cells.Item(74, 4).Value = 44769
$ws.Cells.Item(74, 5).Value = 7
$ws.Cells.Item(74, 6).Value = 100112013
$ws.Cells.Item(74, 7).Value = "Alcachofa"
$ws.Cells.Item(74, 8).Value = "Madrigal"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 300
$ws.Cells.Item(74, 11).Value = 12000
$ws.Cells.Item(74, 12).Value = 12000
$ws.Cells.Item(74, 13).Value = 12000
$ws.Cells.Item(74, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(74, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(74, 16).Value = 300
$ws.Cells.Item(74, 17).Value = 40
$ws.Cells.Item(74, 18).Value = "Hortaliza"
